# Weekly fruit/vegetable price data refresh: the 24 data rows (rows 2-25)
# are reshuffled to a new order (each target row receives the full
# contents of a specific source row from the original layout).
#
# Mapping: targetRow -> sourceRow (1-based spreadsheet row numbers)
$map = @{
    2  = 21
    3  = 19
    4  = 5
    5  = 9
    6  = 13
    7  = 14
    8  = 8
    9  = 25
    10 = 3
    11 = 7
    12 = 10
    13 = 12
    14 = 15
    15 = 6
    16 = 11
    17 = 4
    18 = 24
    19 = 22
    20 = 2
    21 = 23
    22 = 18
    23 = 17
    24 = 16
    25 = 20
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = "A"
$lastCol = "R"
$firstRow = 2
$lastRow = 25

# First, snapshot the full contents of every source row before any writes
# happen, so that overlapping reads/writes in the permutation don't
# clobber data that is still needed. (Number formats stay put per-column,
# e.g. the date format on column D, so only values need to move.)
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rng = $ws.Range("$firstCol$r`:$lastCol$r")
    $snapshot[$r] = $rng.Value2
}

# Now write each target row with the snapshot captured from its mapped
# source row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $map[$r]
    $destRng = $ws.Range("$firstCol$r`:$lastCol$r")
    $destRng.Value = $snapshot[$src]
}
